$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '70.799.02'
$ws.Cells.Item(2, 5).Value = '  -2.02%  '
$ws.Cells.Item(3, 4).Value = '2.551.54'
$ws.Cells.Item(3, 5).Value = '  -5.76%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '578.97'
$ws.Cells.Item(5, 5).Value = '  -3.57%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '170.45'
$ws.Cells.Item(6, 5).Value = '  -3.01%  '
$ws.Cells.Item(7, 5).Value = '  -0.03%  '
$ws.Cells.Item(8, 5).Value = '  -2.68%  '
$ws.Cells.Item(9, 5).Value = '  -1.10%  '
$ws.Cells.Item(10, 4).Value = '2.551.60'
$ws.Cells.Item(10, 5).Value = '  -5.77%  '
$ws.Cells.Item(11, 5).Value = '  -0.27%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.349'
$ws.Cells.Item(12, 5).Value = '  -1.54%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.83'
$ws.Cells.Item(13, 5).Value = '  -3.71%  '
$ws.Cells.Item(14, 4).Value = '3.054.78'
$ws.Cells.Item(14, 5).Value = '  -4.65%  '
$ws.Cells.Item(15, 5).Value = '  -0.75%  '
$ws.Cells.Item(16, 4).Value = '70.664.47'
$ws.Cells.Item(16, 5).Value = '  -1.89%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '25.20'
$ws.Cells.Item(17, 5).Value = '  -4.38%  '
$ws.Cells.Item(18, 4).Value = '2.581.73'
$ws.Cells.Item(18, 5).Value = '  -4.80%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.79'
$ws.Cells.Item(19, 5).Value = '  -3.97%  '
$ws.Cells.Item(20, 5).Value = '  -5.54%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '364.33'
$ws.Cells.Item(21, 5).Value = '  -2.32%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.01'
$ws.Cells.Item(22, 5).Value = '  -4.31%  '
$ws.Cells.Item(23, 5).Value = '  -1.52%  '
$ws.Cells.Item(24, 5).Value = '  -0.16%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '70.05'
$ws.Cells.Item(25, 5).Value = '  -3.31%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '4.15'
$ws.Cells.Item(26, 5).Value = '  -5.28%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.29'
$ws.Cells.Item(27, 5).Value = '  -5.44%  '
$ws.Cells.Item(28, 4).Value = '2.689.36'
$ws.Cells.Item(28, 5).Value = '  -5.55%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.00'
$ws.Cells.Item(29, 5).Value = '  +0.18%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0929'
$ws.Cells.Item(30, 5).Value = '  -6.42%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '7.80'
$ws.Cells.Item(31, 5).Value = '  -4.15%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '484.06'
$ws.Cells.Item(32, 5).Value = '  -4.93%  '
$ws.Cells.Item(33, 5).Value = '  -0.97%  '
$ws.Cells.Item(34, 5).Value = '  -3.22%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.999'
$ws.Cells.Item(35, 5).Value = '  -0.09%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '157.13'
$ws.Cells.Item(36, 5).Value = '  -4.25%  '
$ws.Cells.Item(37, 5).Value = '  +4.61%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '18.79'
$ws.Cells.Item(38, 5).Value = '  -4.64%  '
$ws.Cells.Item(39, 5).Value = '  -1.42%  '
$ws.Cells.Item(40, 5).Value = '  -4.54%  '
$ws.Cells.Item(41, 5).Value = '  -0.07%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.70'
$ws.Cells.Item(42, 5).Value = '  -6.04%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.48'
$ws.Cells.Item(43, 5).Value = '  -3.36%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '4.77'
$ws.Cells.Item(44, 5).Value = '  -5.89%  '
$ws.Cells.Item(45, 5).Value = '  -4.55%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '38.55'
$ws.Cells.Item(46, 5).Value = '  -2.44%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '146.56'
$ws.Cells.Item(47, 5).Value = '  -6.66%  '
$ws.Cells.Item(48, 5).Value = '  -4.97%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.531'
$ws.Cells.Item(49, 5).Value = '  -6.29%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.64'
$ws.Cells.Item(50, 5).Value = '  -7.82%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.595'
$ws.Cells.Item(51, 5).Value = '  -2.17%  '
